$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '45.570.17'
$ws.Range("E2").Value = '  +6.78%  '

$ws.Range("D3").Value = '2.383.95'
$ws.Range("E3").Value = '  +4.70%  '

$ws.Range("E4").Value = '  -1.11%  '

$ws.Range("B5").Value = 'BNB'
$ws.Range("C5").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '318.20'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.79%  '

$ws.Range("B6").Value = 'Solana'
$ws.Range("C6").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '111.39'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +7.66%  '

$ws.Range("E7").Value = '  +2.72%  '

$ws.Range("E8").Value = '  -0.36%  '

$ws.Range("E9").Value = '  +5.74%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.93'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +8.56%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0932'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.90%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.67'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +5.82%  '

$ws.Range("E13").Value = '  +5.37%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.108'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.63%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.75'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +5.24%  '

$ws.Range("D16").Value = '2.742.57'
$ws.Range("E16").Value = '  +4.36%  '

$ws.Range("D17").Value = '2.402.17'
$ws.Range("E17").Value = '  +5.80%  '

$ws.Range("D18").Value = '45.561.53'
$ws.Range("E18").Value = '  +7.07%  '

$ws.Range("E19").Value = '  +5.75%  '

$ws.Range("E20").Value = '  +3.99%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.01'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.77%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '75.11'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.26%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.55'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +4.62%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '269.48'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.87%  '

$ws.Range("E25").Value = '  +7.80%  '

$ws.Range("E26").Value = '  -0.45%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.33'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +6.73%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.55'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +9.41%  '

$ws.Range("E29").Value = '  +0.33%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '23.00'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.94%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '38.71'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +9.06%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0944'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +10.68%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '169.61'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.11%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.05'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +19.20%  '

$ws.Range("E35").Value = '  +2.91%  '

$ws.Range("E36").Value = '  +9.00%  '

$ws.Range("E37").Value = '  +5.46%  '

$ws.Range("E38").Value = '  +12.92%  '

$ws.Range("E39").Value = '  +5.58%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.92'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +5.99%  '

$ws.Range("E41").Value = '  +11.92%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '106.48'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +8.97%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '13.89'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +17.42%  '

$ws.Range("E44").Value = '  +7.07%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '71.39'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +4.25%  '

$ws.Range("E46").Value = '  -0.59%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '118.07'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +7.83%  '

$ws.Range("E48").Value = '  +12.58%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.65'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +21.07%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '79.27'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +5.12%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '9.24'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +7.59%  '
